# edit.ps1 -- "Updated documentation slightly for inlined parameter
# definitions."
#
# Strategy: locate each target paragraph with Range.Find, expand the
# hit to the enclosing paragraph (Range.Expand(4), wdParagraph), and
# replace that whole paragraph's contents in one shot via
# Range.InsertXML with a minimal WordprocessingML package fragment.
# Doing the replacement paragraph-at-a-time (pPr + every run) is what
# lets the new text land split across the exact <w:r> boundaries the
# diff wants (and lets us plant/retire the lone "_GoBack" bookmark)
# without the plain Find/Replace engine silently re-merging runs.
#
# NOTE: `Paragraph.Range` (i.e. `rng.Paragraphs.Item(1).Range`) is not
# reliable in this host -- it does not reflect the real paragraph
# bounds. So paragraph navigation below is done purely in terms of
# Range start/end via Collapse/Move/Expand, never via Paragraph
# objects.

$d = $word.ActiveDocument

function Get-PackageXml([string]$innerXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $innerXml + '</w:p></w:body>' +
           '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Find `needle` in $d.Content and return a Range expanded to the whole
# enclosing paragraph (including the end-of-paragraph mark).
function Get-FullParagraphRange([string]$needle) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find paragraph containing: $needle"
    }
    $rng.Expand(4) | Out-Null   # wdParagraph
    return $rng
}

function Replace-ParagraphXml([string]$needle, [string]$innerXml) {
    $pr = Get-FullParagraphRange $needle
    $pr.InsertXML((Get-PackageXml $innerXml))
}

# ---------------------------------------------------------------------
# 1) Heading: "Parameter management" -> "Parameter " + "definition"
# ---------------------------------------------------------------------
$inner1 = '<w:pPr><w:pStyle w:val="Heading3"/><w:ind w:left="360"/></w:pPr>' +
          '<w:r><w:t xml:space="preserve">Parameter </w:t></w:r>' +
          '<w:r><w:t>definition</w:t></w:r>'
Replace-ParagraphXml "Parameter management" $inner1

# ---------------------------------------------------------------------
# 2) "...Scenario-independent parameters are loaded and derived." ->
#    "...Scenario-independent parameters[_GoBack] are defined or derived."
# ---------------------------------------------------------------------
$inner2 = '<w:pPr><w:ind w:left="720"/></w:pPr>' +
          '<w:r><w:t xml:space="preserve">Parameters defining both the baseline scenario and optional counterfactual scenario are unpacked.  Default values are set for unspecified parameters.  </w:t></w:r>' +
          '<w:r><w:t>Scenario-independent parameters</w:t></w:r>' +
          '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
          '<w:r><w:t xml:space="preserve"> are </w:t></w:r>' +
          '<w:r><w:t>defined or derived.</w:t></w:r>'
Replace-ParagraphXml "Scenario-independent parameters are loaded" $inner2

# ---------------------------------------------------------------------
# 3) Forward-propagation paragraph: drop the old "_GoBack" bookmark and
#    merge "are represented" + " " into a single run (visible text
#    unchanged).
# ---------------------------------------------------------------------
$inner3 = '<w:pPr><w:ind w:left="720"/></w:pPr>' +
          '<w:r><w:t>The population distribution is forward propagated by year using the optima</w:t></w:r>' +
          '<w:r><w:t>l decision values found above.  P</w:t></w:r>' +
          '<w:r><w:t xml:space="preserve">opulation growth </w:t></w:r>' +
          '<w:r><w:t xml:space="preserve">and decline </w:t></w:r>' +
          '<w:r><w:t>due to births</w:t></w:r>' +
          '<w:r><w:t xml:space="preserve">, deaths, and </w:t></w:r>' +
          '<w:r><w:t xml:space="preserve">immigration </w:t></w:r>' +
          '<w:r><w:t xml:space="preserve">flows </w:t></w:r>' +
          '<w:r><w:t xml:space="preserve">are represented </w:t></w:r>' +
          '<w:r><w:t>as explicit distributional changes</w:t></w:r>' +
          '<w:r><w:t xml:space="preserve">.  </w:t></w:r>' +
          '<w:r><w:t xml:space="preserve">For steady states, </w:t></w:r>' +
          '<w:r><w:t>forward propagation continues</w:t></w:r>' +
          '<w:r><w:t xml:space="preserve"> until the normalized distribution over ages satisfies a</w:t></w:r>' +
          '<w:r><w:t xml:space="preserve">n invariance </w:t></w:r>' +
          '<w:r><w:t xml:space="preserve">criterion.  For transition paths, </w:t></w:r>' +
          '<w:r><w:t xml:space="preserve">forward propagation is performed </w:t></w:r>' +
          '<w:r><w:t>across the years</w:t></w:r>' +
          '<w:r><w:t xml:space="preserve"> of the modeling period.</w:t></w:r>' +
          '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
$fwdRng = Get-FullParagraphRange "immigration flows are represented"
$fwdRng.InsertXML((Get-PackageXml $inner3))

# ---------------------------------------------------------------------
# 4) Heading right after the forward-propagation paragraph:
#    "Aggregate" + " generation" -> single run "Aggregate generation"
#    (there are other "Aggregate generation" headings elsewhere in the
#    document, so navigate relative to the paragraph above rather than
#    searching for this text directly).
# ---------------------------------------------------------------------
$headingRng = Get-FullParagraphRange "immigration flows are represented"
$headingRng.Collapse(0) | Out-Null     # wdCollapseEnd
$headingRng.Move(4, 1) | Out-Null      # wdParagraph, forward 1 paragraph
$headingRng.Expand(4) | Out-Null       # wdParagraph
$inner4 = '<w:pPr><w:pStyle w:val="Heading2"/></w:pPr>' +
          '<w:r><w:t>Aggregate generation</w:t></w:r>'
$headingRng.InsertXML((Get-PackageXml $inner4))

# ---------------------------------------------------------------------
# 5) Footer date "...2017-03-21" -> "...2017-04-21", landing the month
#    text split across two runs ("0" and "4") as in the diff.
# ---------------------------------------------------------------------
$footer = $d.Sections.Item(1).Footers.Item(1)
$footerRng = $footer.Range
$rPr5 = '<w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr>'
$inner5 = '<w:pPr><w:pStyle w:val="Footer"/><w:jc w:val="center"/>' + $rPr5 + '</w:pPr>' +
          '<w:r>' + $rPr5 + '<w:t>Pete | 2017</w:t></w:r>' +
          '<w:r>' + $rPr5 + '<w:t>-</w:t></w:r>' +
          '<w:r>' + $rPr5 + '<w:t>0</w:t></w:r>' +
          '<w:r>' + $rPr5 + '<w:t>4</w:t></w:r>' +
          '<w:r>' + $rPr5 + '<w:t>-</w:t></w:r>' +
          '<w:r>' + $rPr5 + '<w:t>21</w:t></w:r>'
$footerRng.InsertXML((Get-PackageXml $inner5))

Write-Output "Done."
